$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    3  = @(0.00006486019690155054, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1, 1.264997119975989)
    4  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    5  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 12.0302756157461)
    6  = @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.145820798638228)
    7  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    8  = @(1.505614041169197, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 3.055818435266709)
    9  = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 0, 22.31973251085698)
    10 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 6.741336633845642)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
